# Insert a new data row at sheet row 170 (pushes existing rows 170:225 down to 171:226)
# and populate it with the new weekly price observation, matching the
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 170:225 down by one row, copying formatting from the row above
# (this is what Excel's own Insert does, and matches style s="2" staying on column D).
$ws.Rows("170:170").Insert()

# Populate the newly inserted row 170 with the new record.
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44627
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = 100112021
$ws.Range("G170").Value = "Ají"
$ws.Range("H170").Value = "Inferno"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 35
$ws.Range("K170").Value = 21000
$ws.Range("L170").Value = 21000
$ws.Range("M170").Value = 21000
$ws.Range("N170").Value = "$/caja 15 kilos"
$ws.Range("O170").Value = "Región Metropolitana"
$ws.Range("P170").Value = 1400
$ws.Range("Q170").Value = 15
$ws.Range("R170").Value = "Hortaliza"
